$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.953047333333333
$ws.Range("H2").Value = 14.859142
$ws.Range("I2").Value = 0.7703204220313993
$ws.Range("J2").Value = 0.7703204220313993
$ws.Range("M2").Value = 3.087329333333333
$ws.Range("N2").Value = 9.261987999999999
$ws.Range("O2").Value = 0.1539049749041678
$ws.Range("P2").Value = 0.1539049749041678
$ws.Range("Q2").Value = 15.29168832158844
$ws.Range("R2").Value = 137.625194894296
$ws.Range("S2").Value = 0.1185561452209104
$ws.Range("T2").Value = 0.1185561452209104

$ws.Range("G3").Value = 4.953047333333333
$ws.Range("H3").Value = 14.859142
$ws.Range("I3").Value = 0.7703204220313993
$ws.Range("J3").Value = 0.7703204220313993
$ws.Range("O3").Value = 0.2832552948356705
$ws.Range("P3").Value = 0.2832552948356705
$ws.Range("Q3").Value = 28.14367558140199
$ws.Range("R3").Value = 253.2930802326179
$ws.Range("S3").Value = 0.2181973382604422
$ws.Range("T3").Value = 0.2181973382604421

$ws.Range("G4").Value = 4.953047333333333
$ws.Range("H4").Value = 14.859142
$ws.Range("I4").Value = 0.7703204220313993
$ws.Range("J4").Value = 0.7703204220313993
$ws.Range("M4").Value = 4.823431
$ws.Range("N4").Value = 14.470293
$ws.Range("O4").Value = 0.2404505470122564
$ws.Range("P4").Value = 0.2404505470122564
$ws.Range("Q4").Value = 23.89068205206733
$ws.Range("R4").Value = 215.016138468606
$ws.Range("S4").Value = 0.1852239668521622
$ws.Range("T4").Value = 0.1852239668521622

$ws.Range("G5").Value = 4.953047333333333
$ws.Range("H5").Value = 14.859142
$ws.Range("I5").Value = 0.7703204220313993
$ws.Range("J5").Value = 0.7703204220313993
$ws.Range("M5").Value = 6.467117666666667
$ws.Range("N5").Value = 19.401353
$ws.Range("O5").Value = 0.3223891832479054
$ws.Range("P5").Value = 0.3223891832479053
$ws.Range("Q5").Value = 32.03193991323622
$ws.Range("R5").Value = 288.287459219126
$ws.Range("S5").Value = 0.2483429716978846
$ws.Range("T5").Value = 0.2483429716978846

$ws.Range("I6").Value = 0.135969508894967
$ws.Range("J6").Value = 0.135969508894967
$ws.Range("M6").Value = 3.087329333333333
$ws.Range("N6").Value = 9.261987999999999
$ws.Range("O6").Value = 0.1539049749041678
$ws.Range("P6").Value = 0.1539049749041678
$ws.Range("Q6").Value = 2.699140892277333
$ws.Range("R6").Value = 24.292268030496
$ws.Range("S6").Value = 0.02092638385421191
$ws.Range("T6").Value = 0.02092638385421191

$ws.Range("I7").Value = 0.135969508894967
$ws.Range("J7").Value = 0.135969508894967
$ws.Range("O7").Value = 0.2832552948356705
$ws.Range("P7").Value = 0.2832552948356705
$ws.Range("S7").Value = 0.0385140833307052
$ws.Range("T7").Value = 0.0385140833307052

$ws.Range("I8").Value = 0.135969508894967
$ws.Range("J8").Value = 0.135969508894967
$ws.Range("M8").Value = 4.823431
$ws.Range("N8").Value = 14.470293
$ws.Range("O8").Value = 0.2404505470122564
$ws.Range("P8").Value = 0.2404505470122564
$ws.Range("Q8").Value = 4.216952079784001
$ws.Range("R8").Value = 37.952568718056
$ws.Range("S8").Value = 0.03269394279078269
$ws.Range("T8").Value = 0.03269394279078268

$ws.Range("I9").Value = 0.135969508894967
$ws.Range("J9").Value = 0.135969508894967
$ws.Range("M9").Value = 6.467117666666667
$ws.Range("N9").Value = 19.401353
$ws.Range("O9").Value = 0.3223891832479054
$ws.Range("P9").Value = 0.3223891832479053
$ws.Range("Q9").Value = 5.653968159730667
$ws.Range("R9").Value = 50.885713437576
$ws.Range("S9").Value = 0.04383509891926722
$ws.Range("T9").Value = 0.04383509891926721

$ws.Range("G10").Value = 0.5382536666666667
$ws.Range("H10").Value = 1.614761
$ws.Range("I10").Value = 0.08371165542397027
$ws.Range("J10").Value = 0.08371165542397027
$ws.Range("M10").Value = 3.087329333333333
$ws.Range("N10").Value = 9.261987999999999
$ws.Range("O10").Value = 0.1539049749041678
$ws.Range("P10").Value = 0.1539049749041678
$ws.Range("Q10").Value = 1.661766333874222
$ws.Range("R10").Value = 14.955897004868
$ws.Range("S10").Value = 0.01288364022721248
$ws.Range("T10").Value = 0.01288364022721248

$ws.Range("G11").Value = 0.5382536666666667
$ws.Range("H11").Value = 1.614761
$ws.Range("I11").Value = 0.08371165542397027
$ws.Range("J11").Value = 0.08371165542397027
$ws.Range("O11").Value = 0.2832552948356705
$ws.Range("P11").Value = 0.2832552948356705
$ws.Range("Q11").Value = 3.058407391591
$ws.Range("R11").Value = 27.525666524319
$ws.Range("S11").Value = 0.02371176963829876
$ws.Range("T11").Value = 0.02371176963829875

$ws.Range("G12").Value = 0.5382536666666667
$ws.Range("H12").Value = 1.614761
$ws.Range("I12").Value = 0.08371165542397027
$ws.Range("J12").Value = 0.08371165542397027
$ws.Range("M12").Value = 4.823431
$ws.Range("N12").Value = 14.470293
$ws.Range("O12").Value = 0.2404505470122564
$ws.Range("P12").Value = 0.2404505470122564
$ws.Range("Q12").Value = 2.596229421663667
$ws.Range("R12").Value = 23.366064794973
$ws.Range("S12").Value = 0.02012851333799517
$ws.Range("T12").Value = 0.02012851333799517

$ws.Range("G13").Value = 0.5382536666666667
$ws.Range("H13").Value = 1.614761
$ws.Range("I13").Value = 0.08371165542397027
$ws.Range("J13").Value = 0.08371165542397027
$ws.Range("M13").Value = 6.467117666666667
$ws.Range("N13").Value = 19.401353
$ws.Range("O13").Value = 0.3223891832479054
$ws.Range("P13").Value = 0.3223891832479053
$ws.Range("Q13").Value = 3.480949796848112
$ws.Range("R13").Value = 31.328548171633
$ws.Range("S13").Value = 0.02698773222046387
$ws.Range("T13").Value = 0.02698773222046386

$ws.Range("G14").Value = 0.06428833333333334
$ws.Range("H14").Value = 0.192865
$ws.Range("I14").Value = 0.009998413649663342
$ws.Range("J14").Value = 0.009998413649663342
$ws.Range("M14").Value = 3.087329333333333
$ws.Range("N14").Value = 9.261987999999999
$ws.Range("O14").Value = 0.1539049749041678
$ws.Range("P14").Value = 0.1539049749041678
$ws.Range("Q14").Value = 0.1984792572911111
$ws.Range("R14").Value = 1.78631331562
$ws.Range("S14").Value = 0.001538805601832925
$ws.Range("T14").Value = 0.001538805601832925

$ws.Range("G15").Value = 0.06428833333333334
$ws.Range("H15").Value = 0.192865
$ws.Range("I15").Value = 0.009998413649663342
$ws.Range("J15").Value = 0.009998413649663342
$ws.Range("O15").Value = 0.2832552948356705
$ws.Range("P15").Value = 0.2832552948356705
$ws.Range("Q15").Value = 0.365292288815
$ws.Range("R15").Value = 3.287630599335
$ws.Range("S15").Value = 0.002832103606224382
$ws.Range("T15").Value = 0.002832103606224382

$ws.Range("G16").Value = 0.06428833333333334
$ws.Range("H16").Value = 0.192865
$ws.Range("I16").Value = 0.009998413649663342
$ws.Range("J16").Value = 0.009998413649663342
$ws.Range("M16").Value = 4.823431
$ws.Range("N16").Value = 14.470293
$ws.Range("O16").Value = 0.2404505470122564
$ws.Range("P16").Value = 0.2404505470122564
$ws.Range("Q16").Value = 0.3100903399383333
$ws.Range("R16").Value = 2.790813059445
$ws.Range("S16").Value = 0.002404124031316362
$ws.Range("T16").Value = 0.002404124031316362

$ws.Range("G17").Value = 0.06428833333333334
$ws.Range("H17").Value = 0.192865
$ws.Range("I17").Value = 0.009998413649663342
$ws.Range("J17").Value = 0.009998413649663342
$ws.Range("M17").Value = 6.467117666666667
$ws.Range("N17").Value = 19.401353
$ws.Range("O17").Value = 0.3223891832479054
$ws.Range("P17").Value = 0.3223891832479053
$ws.Range("Q17").Value = 0.4157602162605556
$ws.Range("R17").Value = 3.741841946345
$ws.Range("S17").Value = 0.003223380410289673
$ws.Range("T17").Value = 0.003223380410289673
